# Implementación del caso de uso "Cambiar alumno de grupo" en el proyecto.
# Actualiza la hoja "Casos de Uso" de la Lista de Tareas de la 5ta Iteración:
#   - La tarea "CU 11 - Cambiar alumno de grupo." (fila 14) registra 5 horas
#     consumidas el Día 3 (columna T) en vez de 4, lo que hace que se
#     recalculen en cascada las columnas "Restante" del resto de días y los
#     totales (AZ14/BA14) de esa fila.
#   - Se actualiza la selección activa de la hoja a la celda T15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Actualizar horas consumidas el Día 3 para "CU 11 - Cambiar alumno de grupo." ---
$ws.Range("T14").Value = 5

# --- Forzar el reordenamiento de las celdas combinadas de la fila de encabezado
#     (Excel las vuelve a escribir al final cuando se tocan) ---
$mergedRanges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($addr in $mergedRanges) {
    $rng = $ws.Range($addr)
    $rng.UnMerge()
    $rng.Merge()
}

# --- Actualizar la selección activa de la hoja a T15 ---
$ws.Range("T15").Select()

$wb.Save()
